# "Completed Reorganization and Ch 3 edits"
#
# The document is a Word master document whose front-matter / chapter
# sections each carry a cached " PAGE " field result in their
# header/footer (the last page number Word computed for that section
# the last time the whole master + subdocuments were laid out). The
# chapter reorg shifted those cached page numbers:
#
#   - List of Abbreviations section footer: 12 -> 22
#   - Chapter 1 section footer:            28 -> 23
#   - Chapter 2 section header:            28 -> 15
#
# Each value lives inside a PAGE field's cached result run
# (<w:t> inside the run after the "separate" fldChar), so we scope a
# Find/Replace to the specific header/footer range that owns each
# field, matching the whole (one-character-run) field result text so we
# only ever touch that single cached digit run.

$d = $word.ActiveDocument

# NOTE: this interpreter's function-call binding for *named* parameters
# (e.g. "-Range $r") does not reliably hand COM range objects through to
# the callee, so the helper below takes plain positional parameters.
function Set-CachedPageField {
    param($Range, $OldText, $NewText)

    $Range.Find.Execute($OldText, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $NewText, 2) | Out-Null
}

# Section 9  -> abbrev.docx (List of Abbreviations) footer: 12 -> 22
Set-CachedPageField $d.Sections.Item(9).Footers.Item(1).Range "12" "22"

# Section 10 -> Chapter1.docx footer: 28 -> 23
Set-CachedPageField $d.Sections.Item(10).Footers.Item(1).Range "28" "23"

# Section 11 -> Chapter2.docx header: 28 -> 15
Set-CachedPageField $d.Sections.Item(11).Headers.Item(1).Range "28" "15"
